$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C (so everything from C onward shifts right).
$ws.Range("C1").EntireColumn.Insert()

# Update header row text.
$ws.Range("B1").Value = "RUBRO TEMPORAL"
$ws.Range("C1").Value = "RUBRO PERMANENTE"

# Restore/confirm the selection shown by Excel after the edit.
$ws.Range("B2").Select()
